$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 105.9786

$ws.Range("C3").Value = 299.0947
$ws.Range("C4").Value = 525.1708

$ws.Range("C5").Value = 587.8281
$ws.Range("C5").Font.Bold = $false

$ws.Range("C6").Value = 1269.227
$ws.Range("C7").Value = 1190.797
$ws.Range("C8").Value = 2919
$ws.Range("C9").Value = 2178.051

$ws.Range("G12").Select()
